$wb = $excel.ActiveWorkbook

# --- loandetails sheet: refresh the group/mobile test rows ---
$loandetails = $wb.Worksheets.Item("loandetails")

# Row 2 -> group005
$loandetails.Range("A2").Value = "'6000010005"
$loandetails.Range("B2").Value = "testAutomationg005"
$loandetails.Range("C2").Value = "testAutomationg005"
$loandetails.Range("D2").Value = "testAutomationg005"

# Row 3 -> group001
$loandetails.Range("A3").Value = "'6000010000"
$loandetails.Range("B3").Value = "testAutomationg001"
$loandetails.Range("C3").Value = "testAutomationg001"
$loandetails.Range("D3").Value = "testAutomationg001"

# Row 4 -> group003 (content unchanged, but rewritten same as before)
$loandetails.Range("A4").Value = "'6000010002"
$loandetails.Range("B4").Value = "testAutomationg003"
$loandetails.Range("C4").Value = "testAutomationg003"
$loandetails.Range("D4").Value = "testAutomationg003"

# Row 5 -> group002, now also gets a value in column A
$loandetails.Range("A5").Value = "'6000010001"
$loandetails.Range("B5").Value = "testAutomationg002"
$loandetails.Range("C5").Value = "testAutomationg002"
$loandetails.Range("D5").Value = "testAutomationg002"

# Row 6 -> group006, now also gets a value in column A
$loandetails.Range("A6").Value = "'6000010006"
$loandetails.Range("B6").Value = "testAutomationg006"
$loandetails.Range("C6").Value = "testAutomationg006"
$loandetails.Range("D6").Value = "testAutomationg006"

# --- sheet selection / active-tab bookkeeping ---
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Activate() | Out-Null
$sheet1.Range("B3").Select() | Out-Null

$loandetails.Activate() | Out-Null
$loandetails.Range("A10").Select() | Out-Null
